$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("F24").Select()
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
